$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (columns A:F) as described by the diff.
$data = @(
    @(45959, 8119.74184970922, 7622.16034200808, 10844, 3328.43745987902, 4.44157507862891),
    @(45960, 7630.16423365932, 7006.4744210815, 3620, 3717.60068615447, 296.003129468165),
    @(45961, 7043.69124402393, 5811.11699649386, 3620, 3317.85730607899, 229.540595940535),
    @(45962, 2165.99896552164, 3487.59902092829, 11404, 3351.97858663452, -190.18426635155),
    @(45963, 2123.07986091365, 3749.81266593743, 11404, 3310.05866079226, -181.005361386263),
    @(45964, 8571.71877057272, 7322.85732464797, 11404, 4686.02515008376, 25.2034364471553),
    @(45965, 8571.71877057272, 7666.32173668027, 11404, 4686.02515008376, 39.5144536151679),
    @(45966, 8571.71877057272, 7818.38260030498, 11404, 4686.02515008376, 45.8503229328641),
    @(45967, 8580.99785695336, 7996.45957611782, 11404, 4688.11424399391, 53.3572425046558),
    @(45968, 8585.92244945998, 7205.80789219939, 11404, 4719.63766685978, 21.7268982941317),
    @(45969, 3147.97579797074, 4895.47337290744, 11404, 4404.41280498406, -87.6714092545208),
    @(45970, 3012.17181266658, 4815.45056894055, 11404, 4396.2731671762, -91.3448443284688),
    @(45971, 9340.71276866398, 8144.87357607735, 11404, 5362.50436831596, 87.6407476830548),
    @(45972, 9340.71407063505, 8641.81501444452, 11404, 5362.50813242099, 108.346797786063)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
